$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Grundmaße")
$ws2 = $wb.Worksheets.Item("Rohre")

# --- Sheet "Grundmaße": title text changed to be PCB-specific ---
$ws1.Range("A1").Value = "Grundsätzliche Daten PCB"

# --- Sheet "Rohre": widen column A to fit the new longer labels ---
$ws2.Columns.Item(1).ColumnWidth = 20.5

# New row about ALU advantages
$ws2.Range("A17").Value = "Vorteile ALU"

# New row about the PVC supplier (keep shared-string append order: K20 then A20)
$ws2.Range("K20").Value = "Versandkosten 4,90€ pro Bestellung"
$ws2.Range("A20").Value = "Lieferant GWT für PVC"

# H20 gets an "Amazon"-style hyperlink like the other rows above
$ws2.Hyperlinks.Add($ws2.Range("H20"), "https://www.gwt-kunststoffe.de/pvc-rundstab", "", "", "Amazon")
$ws2.Range("H4").Copy()
$ws2.Range("H20").PasteSpecial(-4122)

# Restore selections as left by the author
$ws1.Activate()
$ws1.Range("A2").Select()
$ws2.Activate()
$ws2.Range("G23").Select()
